$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The two 8:xx time-slot labels (shared strings used by C2/C3) were shifted
# by 15 minutes: "8:25-8:30" -> "8:40-8:45" and "8:30-8:35" -> "8:45-8:50".
$ws.Range("C2").Value = "8:40-8:45"
$ws.Range("C3").Value = "8:45-8:50"

# Expand the saved selection from the single cell C10 to the C9:C10 block
# (active cell stays anchored at C10).
$ws.Range("C10").Select()
$ws.Range("C9:C10").Select()
